# Slide 8, "Content Placeholder 2": shrink the three bullet paragraphs from
# 26pt to 24pt. In real PowerPoint, editing the font size causes the
# "shrink text on overflow" autofit to recompute, and since the text now
# fits at 24pt the <a:normAutofit lnSpcReduction="10000"/> no longer needs
# its line-spacing reduction, so PowerPoint drops the attribute, leaving a
# bare <a:normAutofit/>.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text -ne "") {
        $para.Font.Size = 24
    }
}

# Re-assert the "shrink text on overflow" autofit so the host recomputes
# the bodyPr, clearing the stale lnSpcReduction="10000" left over from the
# 26pt layout.
$tf.AutoSize = 2
